$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "30.546.22"
$ws.Cells.Item(2, 5).Value = "  -0.13%  "
Set-TextValue 3 4 "1.917.53"
$ws.Cells.Item(3, 5).Value = "  -0.45%  "
Set-TextValue 4 4 "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
Set-TextValue 5 4 "244.75"
$ws.Cells.Item(5, 5).Value = "  -0.80%  "
$ws.Cells.Item(6, 5).Value = "  +0.05%  "
Set-TextValue 7 4 "0.4812"
$ws.Cells.Item(7, 5).Value = "  +1.57%  "
Set-TextValue 8 4 "0.2897"
$ws.Cells.Item(8, 5).Value = "  -0.65%  "
Set-TextValue 9 4 "0.06704"
$ws.Cells.Item(9, 5).Value = "  -1.47%  "
Set-TextValue 10 4 "111.24"
$ws.Cells.Item(10, 5).Value = "  +5.11%  "
Set-TextValue 11 4 "18.93"
$ws.Cells.Item(11, 5).Value = "  +2.91%  "
Set-TextValue 12 4 "1.919.65"
$ws.Cells.Item(12, 5).Value = "  -0.25%  "
Set-TextValue 13 4 "0.07560"
$ws.Cells.Item(13, 5).Value = "  -2.25%  "
Set-TextValue 14 4 "5.282"
$ws.Cells.Item(14, 5).Value = "  -1.31%  "
Set-TextValue 15 4 "0.6672"
$ws.Cells.Item(15, 5).Value = "  -0.93%  "
Set-TextValue 16 4 "295.35"
$ws.Cells.Item(16, 5).Value = "  +2.18%  "
Set-TextValue 17 4 "30.547.96"
$ws.Cells.Item(17, 5).Value = "  -0.23%  "
$ws.Cells.Item(18, 2).Value = "Dai"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 18 4 "1.001"
$ws.Cells.Item(18, 5).Value = "  +0.02%  "
$ws.Cells.Item(19, 2).Value = "ShibaInu"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue 19 4 "0.000007599"
$ws.Cells.Item(19, 5).Value = "  -0.81%  "
Set-TextValue 20 4 "12.96"
$ws.Cells.Item(20, 5).Value = "  -0.68%  "
Set-TextValue 21 4 "5.539"
$ws.Cells.Item(21, 5).Value = "  +0.67%  "
Set-TextValue 22 4 "2.171.00"
$ws.Cells.Item(22, 5).Value = "  +0.17%  "
Set-TextValue 23 4 "1.002"
$ws.Cells.Item(23, 5).Value = "  +0.08%  "
Set-TextValue 24 4 "6.439"
$ws.Cells.Item(24, 5).Value = "  +2.43%  "
$ws.Cells.Item(25, 5).Value = "  +0.42%  "
Set-TextValue 26 4 "165.17"
$ws.Cells.Item(26, 5).Value = "  -2.23%  "
Set-TextValue 27 4 "20.22"
$ws.Cells.Item(27, 5).Value = "  -2.29%  "
$ws.Cells.Item(28, 5).Value = "  -1.97%  "
$ws.Cells.Item(29, 5).Value = "  -2.35%  "
Set-TextValue 30 4 "1.432"
$ws.Cells.Item(30, 5).Value = "  +5.19%  "
Set-TextValue 31 4 "4.138"
$ws.Cells.Item(31, 5).Value = "  -0.97%  "
Set-TextValue 32 4 "4.061"
$ws.Cells.Item(32, 5).Value = "  -0.30%  "
Set-TextValue 33 4 "0.04997"
$ws.Cells.Item(33, 5).Value = "  -1.51%  "
Set-TextValue 34 4 "0.7383"
$ws.Cells.Item(34, 5).Value = "  -0.43%  "
Set-TextValue 35 4 "1.134"
$ws.Cells.Item(35, 5).Value = "  -1.79%  "
Set-TextValue 36 4 "0.9997"
$ws.Cells.Item(36, 5).Value = "  -0.03%  "
Set-TextValue 37 4 "2.722"
$ws.Cells.Item(37, 5).Value = "  -0.56%  "
Set-TextValue 39 4 "2.681"
$ws.Cells.Item(39, 5).Value = "  -0.52%  "
Set-TextValue 40 4 "110.64"
$ws.Cells.Item(40, 5).Value = "  -0.60%  "
Set-TextValue 41 4 "2.009"
$ws.Cells.Item(41, 5).Value = "  -2.59%  "
Set-TextValue 42 4 "0.4415"
$ws.Cells.Item(42, 5).Value = "  -1.88%  "
Set-TextValue 43 4 "0.8652"
$ws.Cells.Item(43, 5).Value = "  -1.20%  "
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 44 4 "5.831"
$ws.Cells.Item(44, 5).Value = "  -1.66%  "
$ws.Cells.Item(45, 2).Value = "Aave"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue 45 4 "70.33"
$ws.Cells.Item(45, 5).Value = "  +3.86%  "
$ws.Cells.Item(46, 5).Value = "  +0.04%  "
Set-TextValue 47 4 "7.206"
$ws.Cells.Item(47, 5).Value = "  -1.11%  "
Set-TextValue 48 4 "48.56"
$ws.Cells.Item(48, 5).Value = "  +2.55%  "
Set-TextValue 49 4 "9.173"
$ws.Cells.Item(49, 5).Value = "  -2.12%  "
$ws.Cells.Item(50, 5).Value = "  -0.57%  "
$ws.Cells.Item(51, 5).Value = "  -0.96%  "
